$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 4584.25
$ws.Cells.Item(43, 9).Value = 3404
$ws.Cells.Item(43, 10).Value = 5292.4
$ws.Cells.Item(43, 11).Value = 3404
$ws.Cells.Item(43, 12).Value = 5292.4
$ws.Cells.Item(43, 13).Value = -3335
$ws.Cells.Item(43, 14).Value = -5430.4

$ws.Cells.Item(62, 8).Value = 4395.543
$ws.Cells.Item(62, 9).Value = 3948.3333
$ws.Cells.Item(62, 10).Value = 7078.8
$ws.Cells.Item(62, 11).Value = 3948.3333
$ws.Cells.Item(62, 12).Value = 7078.8
$ws.Cells.Item(62, 13).Value = -3324.3333
$ws.Cells.Item(62, 14).Value = -8326.799999999999

$ws.Cells.Item(65, 8).Value = 4395.543
$ws.Cells.Item(65, 9).Value = 3948.3333
$ws.Cells.Item(65, 10).Value = 7078.8
$ws.Cells.Item(65, 11).Value = 19741.6665
$ws.Cells.Item(65, 12).Value = 35394
$ws.Cells.Item(65, 13).Value = -16621.6665
$ws.Cells.Item(65, 14).Value = -41634

$ws.Cells.Item(70, 8).Value = 52086.363
$ws.Cells.Item(70, 10).Value = 82884.766
$ws.Cells.Item(70, 12).Value = 248654.298
$ws.Cells.Item(70, 14).Value = -249194.298

$ws.Cells.Item(73, 8).Value = 52086.363
$ws.Cells.Item(73, 10).Value = 82884.766
$ws.Cells.Item(73, 12).Value = 248654.298
$ws.Cells.Item(73, 14).Value = -250526.298

$ws.Cells.Item(125, 8).Value = 102343.664
$ws.Cells.Item(125, 10).Value = 6000
$ws.Cells.Item(125, 12).Value = 54000
$ws.Cells.Item(125, 14).Value = -58920

$ws.Cells.Item(132, 8).Value = 32443.39
$ws.Cells.Item(132, 9).Value = 33691.047
$ws.Cells.Item(132, 11).Value = 101073.141
$ws.Cells.Item(132, 13).Value = -98543.141

$ws.Cells.Item(135, 8).Value = 1257.2354
$ws.Cells.Item(135, 9).Value = 562.3570999999999
$ws.Cells.Item(135, 11).Value = 5061.2139
$ws.Cells.Item(135, 13).Value = -2526.2139

$ws.Cells.Item(136, 8).Value = 0
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 14).ClearContents()

$ws.Cells.Item(137, 8).Value = 15479.782
$ws.Cells.Item(137, 9).Value = 19937.705
$ws.Cells.Item(137, 11).Value = 59813.11500000001
$ws.Cells.Item(137, 13).Value = -57263.11500000001

$ws.Cells.Item(140, 8).Value = 95991
$ws.Cells.Item(140, 10).Value = 95991
$ws.Cells.Item(140, 12).Value = 95991
$ws.Cells.Item(140, 14).Value = -106351

$ws.Cells.Item(141, 8).Value = 1450.6666
$ws.Cells.Item(141, 9).Value = 1399.8
$ws.Cells.Item(141, 11).Value = 4199.4
$ws.Cells.Item(141, 13).Value = 980.6000000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 11459.111
$ws.Cells.Item(5, 9).Value = 12843.375
$ws.Cells.Item(5, 11).Value = 12843.375
$ws.Cells.Item(5, 13).Value = -12731.375

$ws.Cells.Item(32, 8).Value = 20929.648
$ws.Cells.Item(32, 9).Value = 22395.62
$ws.Cells.Item(32, 11).Value = 22395.62
$ws.Cells.Item(32, 13).Value = -22108.62

$ws.Cells.Item(45, 8).Value = 4050.5625
$ws.Cells.Item(45, 9).Value = 2855.2222
$ws.Cells.Item(45, 11).Value = 2855.2222
$ws.Cells.Item(45, 13).Value = -2478.2222

$ws.Cells.Item(74, 8).Value = 278193.2
$ws.Cells.Item(74, 9).Value = 334097.06
$ws.Cells.Item(74, 10).Value = 26625.75
$ws.Cells.Item(74, 11).Value = 334097.06
$ws.Cells.Item(74, 12).Value = 26625.75
$ws.Cells.Item(74, 13).Value = -333223.06
$ws.Cells.Item(74, 14).Value = -28373.75

$ws.Cells.Item(77, 8).Value = 278193.2
$ws.Cells.Item(77, 9).Value = 334097.06
$ws.Cells.Item(77, 10).Value = 26625.75
$ws.Cells.Item(77, 11).Value = 1670485.3
$ws.Cells.Item(77, 12).Value = 133128.75
$ws.Cells.Item(77, 13).Value = -1666117.3
$ws.Cells.Item(77, 14).Value = -141864.75

$ws.Cells.Item(139, 8).Value = 132999.6
$ws.Cells.Item(139, 10).Value = 132999.6
$ws.Cells.Item(139, 12).Value = 132999.6
$ws.Cells.Item(139, 14).Value = -143279.6

$ws.Cells.Item(141, 8).Value = 92719.25
$ws.Cells.Item(141, 10).Value = 92719.25
$ws.Cells.Item(141, 12).Value = 92719.25
$ws.Cells.Item(141, 14).Value = -103079.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 11459.111
$ws.Cells.Item(4, 9).Value = 12843.375
$ws.Cells.Item(4, 11).Value = 12843.375
$ws.Cells.Item(4, 13).Value = -12728.375

$ws.Cells.Item(20, 8).Value = 23590.5
$ws.Cells.Item(20, 9).Value = 29473.455
$ws.Cells.Item(20, 10).Value = 2019.6666
$ws.Cells.Item(20, 11).Value = 29473.455
$ws.Cells.Item(20, 12).Value = 2019.6666
$ws.Cells.Item(20, 13).Value = -29226.455
$ws.Cells.Item(20, 14).Value = -2513.6666

$ws.Cells.Item(139, 8).Value = 86930.664
$ws.Cells.Item(139, 9).Value = 65709
$ws.Cells.Item(139, 10).Value = 97541.5
$ws.Cells.Item(139, 11).Value = 65709
$ws.Cells.Item(139, 12).Value = 97541.5
$ws.Cells.Item(139, 13).Value = -60569
$ws.Cells.Item(139, 14).Value = -107821.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4350022.5
$ws.Cells.Item(31, 9).Value = 6667721
$ws.Cells.Item(31, 10).Value = 4337.625
$ws.Cells.Item(31, 11).Value = 6667721
$ws.Cells.Item(31, 12).Value = 4337.625
$ws.Cells.Item(31, 13).Value = -6667426
$ws.Cells.Item(31, 14).Value = -4927.625

$ws.Cells.Item(34, 8).Value = 4350022.5
$ws.Cells.Item(34, 9).Value = 6667721
$ws.Cells.Item(34, 10).Value = 4337.625
$ws.Cells.Item(34, 11).Value = 6667721
$ws.Cells.Item(34, 12).Value = 4337.625
$ws.Cells.Item(34, 13).Value = -6667519
$ws.Cells.Item(34, 14).Value = -4741.625

$ws.Cells.Item(94, 8).Value = 1083.9166
$ws.Cells.Item(94, 9).Value = 940.2222
$ws.Cells.Item(94, 10).Value = 1170.1333
$ws.Cells.Item(94, 11).Value = 940.2222
$ws.Cells.Item(94, 12).Value = 1170.1333
$ws.Cells.Item(94, 13).Value = -489.2222
$ws.Cells.Item(94, 14).Value = -2072.1333

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 44194820
$ws.Cells.Item(4, 9).Value = 63296590
$ws.Cells.Item(4, 11).Value = 189889770
$ws.Cells.Item(4, 13).Value = -189889658

$ws.Cells.Item(5, 8).Value = 1050.5294
$ws.Cells.Item(5, 9).Value = 947
$ws.Cells.Item(5, 11).Value = 2841
$ws.Cells.Item(5, 13).Value = -2729

$ws.Cells.Item(116, 8).Value = 2643.6667
$ws.Cells.Item(116, 9).Value = 2999.5
$ws.Cells.Item(116, 10).Value = 1932
$ws.Cells.Item(116, 11).Value = 8998.5
$ws.Cells.Item(116, 12).Value = 5796
$ws.Cells.Item(116, 13).Value = -5556.5
$ws.Cells.Item(116, 14).Value = -12680

$ws.Cells.Item(122, 8).Value = 966.5833
$ws.Cells.Item(122, 9).Value = 338.6
$ws.Cells.Item(122, 11).Value = 3047.4
$ws.Cells.Item(122, 13).Value = -597.4000000000001

$ws.Cells.Item(129, 8).Value = 2420
$ws.Cells.Item(129, 9).Value = 1769.6923
$ws.Cells.Item(129, 10).Value = 3265.4
$ws.Cells.Item(129, 11).Value = 5309.0769
$ws.Cells.Item(129, 12).Value = 9796.200000000001
$ws.Cells.Item(129, 13).Value = -309.0769
$ws.Cells.Item(129, 14).Value = -19796.2

$ws.Cells.Item(132, 8).Value = 1420.1428
$ws.Cells.Item(132, 9).Value = 1172
$ws.Cells.Item(132, 11).Value = 10548
$ws.Cells.Item(132, 13).Value = -8018

$ws.Cells.Item(133, 8).Value = 2734.5833
$ws.Cells.Item(133, 9).Value = 2734.5833
$ws.Cells.Item(133, 11).Value = 8203.749899999999
$ws.Cells.Item(133, 13).Value = -3143.749899999999

$ws.Cells.Item(134, 8).Value = 892.8182
$ws.Cells.Item(134, 9).Value = 892.8182
$ws.Cells.Item(134, 11).Value = 2678.4546
$ws.Cells.Item(134, 13).Value = 2391.5454

$ws.Cells.Item(135, 8).Value = 1050.5294
$ws.Cells.Item(135, 9).Value = 947
$ws.Cells.Item(135, 11).Value = 8523
$ws.Cells.Item(135, 13).Value = -5988

$ws.Cells.Item(136, 8).Value = 1846.875
$ws.Cells.Item(136, 9).Value = 1846.875
$ws.Cells.Item(136, 11).Value = 5540.625
$ws.Cells.Item(136, 13).Value = -440.625

$ws.Cells.Item(137, 8).Value = 3525.389
$ws.Cells.Item(137, 9).Value = 2847.2222
$ws.Cells.Item(137, 10).Value = 4203.5557
$ws.Cells.Item(137, 11).Value = 8541.6666
$ws.Cells.Item(137, 12).Value = 12610.6671
$ws.Cells.Item(137, 13).Value = -3441.6666
$ws.Cells.Item(137, 14).Value = -22810.6671

$ws.Cells.Item(138, 8).Value = 6183.636
$ws.Cells.Item(138, 9).Value = 1498.375
$ws.Cells.Item(138, 11).Value = 4495.125
$ws.Cells.Item(138, 13).Value = 644.875

$ws.Cells.Item(140, 8).Value = 5064.5557
$ws.Cells.Item(140, 9).Value = 5072.625
$ws.Cells.Item(140, 11).Value = 15217.875
$ws.Cells.Item(140, 13).Value = -10037.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(29, 8).Value = 15428.929
$ws.Cells.Item(29, 9).Value = 6601
$ws.Cells.Item(29, 10).Value = 20333.334
$ws.Cells.Item(29, 11).Value = 6601
$ws.Cells.Item(29, 12).Value = 20333.334
$ws.Cells.Item(29, 13).Value = -6311
$ws.Cells.Item(29, 14).Value = -20913.334

$ws.Cells.Item(80, 8).Value = 12431.75
$ws.Cells.Item(80, 9).Value = 6126.25
$ws.Cells.Item(80, 10).Value = 18737.25
$ws.Cells.Item(80, 11).Value = 6126.25
$ws.Cells.Item(80, 12).Value = 18737.25
$ws.Cells.Item(80, 13).Value = -5128.25
$ws.Cells.Item(80, 14).Value = -20733.25

$ws.Cells.Item(83, 8).Value = 12431.75
$ws.Cells.Item(83, 9).Value = 6126.25
$ws.Cells.Item(83, 10).Value = 18737.25
$ws.Cells.Item(83, 11).Value = 30631.25
$ws.Cells.Item(83, 12).Value = 93686.25
$ws.Cells.Item(83, 13).Value = -25639.25
$ws.Cells.Item(83, 14).Value = -103670.25

$ws.Cells.Item(122, 8).Value = 3282.1052
$ws.Cells.Item(122, 9).Value = 3350.1333
$ws.Cells.Item(122, 11).Value = 10050.3999
$ws.Cells.Item(122, 13).Value = -7600.3999

$ws.Cells.Item(132, 8).Value = 3086.875
$ws.Cells.Item(132, 9).Value = 2959.3333
$ws.Cells.Item(132, 10).Value = 5000
$ws.Cells.Item(132, 11).Value = 8877.999899999999
$ws.Cells.Item(132, 12).Value = 15000
$ws.Cells.Item(132, 13).Value = -6347.999899999999
$ws.Cells.Item(132, 14).Value = -20060

$ws.Cells.Item(139, 8).Value = 112442.25
$ws.Cells.Item(139, 10).Value = 112442.25
$ws.Cells.Item(139, 12).Value = 112442.25
$ws.Cells.Item(139, 14).Value = -122722.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 4254.9287
$ws.Cells.Item(46, 9).Value = 884.1667
$ws.Cells.Item(46, 10).Value = 6783
$ws.Cells.Item(46, 11).Value = 884.1667
$ws.Cells.Item(46, 12).Value = 6783
$ws.Cells.Item(46, 13).Value = -696.1667
$ws.Cells.Item(46, 14).Value = -7159

$ws.Cells.Item(82, 8).Value = 1287.1765
$ws.Cells.Item(82, 9).Value = 1278.5714
$ws.Cells.Item(82, 10).Value = 1293.2
$ws.Cells.Item(82, 11).Value = 1278.5714
$ws.Cells.Item(82, 12).Value = 1293.2
$ws.Cells.Item(82, 13).Value = -917.5714
$ws.Cells.Item(82, 14).Value = -2015.2

$ws.Cells.Item(85, 8).Value = 1287.1765
$ws.Cells.Item(85, 9).Value = 1278.5714
$ws.Cells.Item(85, 10).Value = 1293.2
$ws.Cells.Item(85, 11).Value = 1278.5714
$ws.Cells.Item(85, 12).Value = 1293.2
$ws.Cells.Item(85, 13).Value = -30.57140000000004
$ws.Cells.Item(85, 14).Value = -3789.2

$ws.Cells.Item(136, 8).Value = 4010.8635
$ws.Cells.Item(136, 9).Value = 3594.842
$ws.Cells.Item(136, 11).Value = 10784.526
$ws.Cells.Item(136, 13).Value = -8234.526

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(49, 8).Value = 12450
$ws.Cells.Item(49, 9).Value = 4500
$ws.Cells.Item(49, 10).Value = 13333.333
$ws.Cells.Item(49, 11).Value = 4500
$ws.Cells.Item(49, 12).Value = 13333.333
$ws.Cells.Item(49, 13).Value = -4270
$ws.Cells.Item(49, 14).Value = -13793.333

$ws.Cells.Item(132, 8).Value = 2030.8096
$ws.Cells.Item(132, 9).Value = 1597
$ws.Cells.Item(132, 10).Value = 3874.5
$ws.Cells.Item(132, 11).Value = 4791
$ws.Cells.Item(132, 12).Value = 11623.5
$ws.Cells.Item(132, 13).Value = -2261
$ws.Cells.Item(132, 14).Value = -16683.5

$ws.Cells.Item(136, 8).Value = 32108.945
$ws.Cells.Item(136, 9).Value = 33497.766
$ws.Cells.Item(136, 10).Value = 8499
$ws.Cells.Item(136, 11).Value = 100493.298
$ws.Cells.Item(136, 12).Value = 25497
$ws.Cells.Item(136, 13).Value = -97943.29800000001
$ws.Cells.Item(136, 14).Value = -30597
